# Weekly driver report update for 2025-04-19
# Rows 12-17 of the "Good Drivers" table get re-ranked and their
# client counts refreshed; two newly-promoted drivers (21.60.2.1 and
# 22.50.1.1) have not yet accumulated a "Driver Vintage" date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
$ws.Cells.Item(12, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Cells.Item(12, 2).Value = 56018
$ws.Cells.Item(12, 4).Value = 100
$ws.Cells.Item(12, 5).ClearContents()

# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Cells.Item(13, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Cells.Item(13, 2).Value = 34244
$ws.Cells.Item(13, 4).Value = 100
$ws.Cells.Item(13, 5).ClearContents()

# Row 14: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Cells.Item(14, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Cells.Item(14, 2).Value = 442178
$ws.Cells.Item(14, 4).Value = 99.90000000000001
# Leading apostrophe forces the date-shaped string to stay plain text
# instead of being auto-converted into a date serial value.
$ws.Cells.Item(14, 5).Value = "'2024-11-10"

# Row 15: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Cells.Item(15, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Cells.Item(15, 2).Value = 77849
$ws.Cells.Item(15, 4).Value = 99.90000000000001
$ws.Cells.Item(15, 5).Value = "'2021-08-18"

# Row 16: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Cells.Item(16, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Cells.Item(16, 2).Value = 59673
$ws.Cells.Item(16, 4).Value = 100
$ws.Cells.Item(16, 5).Value = "'2020-08-05"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Cells.Item(17, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Cells.Item(17, 2).Value = 113652
$ws.Cells.Item(17, 4).Value = 100
$ws.Cells.Item(17, 5).Value = "'2019-12-14"
